# Re-apply the "wieder auskommentiert" change: the generator now produces
# values around Int32.MaxValue (2147483600 .. 2147483647) for rows 2-49,
# while row 50 keeps the original Int32.MinValue edge case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$onesWords    = @("", "eins", "zwei", "drei", "vier", "fuenf", "sechs", "sieben", "acht", "neun",
                   "zehn", "elf", "zwoelf", "dreizehn", "vierzehn", "fuenfzehn", "sechzehn", "siebzehn", "achtzehn", "neunzehn")
$onesWordsAnd = @("", "ein", "zwei", "drei", "vier", "fuenf", "sechs", "sieben", "acht", "neun")
$tensWords    = @{ 20 = "zwanzig"; 30 = "dreissig"; 40 = "vierzig"; 50 = "fuenfzig";
                    60 = "sechzig"; 70 = "siebzig"; 80 = "achtzig"; 90 = "neunzig" }

function Spell-0-99([int]$n) {
    if ($n -eq 0) { return "" }
    if ($n -lt 20) { return $onesWords[$n] }
    $t = [int]([math]::Floor($n / 10)) * 10
    $u = $n % 10
    if ($u -eq 0) { return $tensWords[$t] }
    return $onesWordsAnd[$u] + "und" + $tensWords[$t]
}

$prefix = "zweimilliardeneinhundertsiebenundvierzigmillionenvierhundertdreiundachtzigtausendsechshundert"

# Rows 2..49: values Int32.MaxValue-47 .. Int32.MaxValue, positive, no "minus"
for ($row = 2; $row -le 49; $row++) {
    $n = $row - 2
    $value = 2147483600 + $n
    $text = " " + $prefix + (Spell-0-99 $n)

    $ws.Cells.Item($row, 1).Value = $value
    $ws.Cells.Item($row, 2).Value = $text
    $ws.Cells.Item($row, 3).Value = $text
}

# Row 50: back to the Int32.MinValue edge case (commented back in)
$minText = " minus " + $prefix + (Spell-0-99 48)
$ws.Cells.Item(50, 1).Value = -2147483648
$ws.Cells.Item(50, 2).Value = $minText
$ws.Cells.Item(50, 3).Value = $minText
